$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 64
$ws.Range("A64").Value = 1339590450999557888
$v = @'
Dua Ribu Dua Pulu-h
Kuliah tura turu
Ujug-ujug tahun baru
#COVID19indonesia #kuliahonline #belajardaring #kuliahdaring #2020 #2021 #Online #OnlineClasses #rebahan #seninsemangat #selasa #rabu #KamisManis #JumatBerkah #sabtu #minggu
'@
$ws.Range("B64").Value = $v
$v = @'
allifi_nasihah
'@
$ws.Range("C64").Value = $v
$v = @'
Thu Dec 17 15:17:06 +0000 2020
'@
$ws.Range("D64").Value = $v
$ws.Rows.Item(64).AutoFit() | Out-Null

# Row 65
$ws.Range("A65").Value = 1339550525361017088
$v = @'
*JAPRI DOSEN*
Mhsw : Mohon maaf pak, apakah boleh untuk pengumpulan UTS nya hari senin? Karena ada beberapa teman kita yg posisinya sedang tidak dikota, dan terhalang untuk mengumpulkan langsung .
Dosen : TERSERAH SUDAH
monggo rek piye masamu?
#kuliahonline #deritadaring
'@
$ws.Range("B65").Value = $v
$v = @'
panggilajasinyo
'@
$ws.Range("C65").Value = $v
$v = @'
Thu Dec 17 12:38:27 +0000 2020
'@
$ws.Range("D65").Value = $v
$ws.Rows.Item(65).AutoFit() | Out-Null

# Row 66
$ws.Range("A66").Value = 1339511005051686912
$v = @'
Tahun depan kuliah online, bantu cari alasan dong supaya bisa balik. Bosan di rumah uy, wkwk 😅😅 #kuliahonline #daring #kuliahdaring
'@
$ws.Range("B66").Value = $v
$v = @'
danielbintangc3
'@
$ws.Range("C66").Value = $v
$v = @'
Thu Dec 17 10:01:25 +0000 2020
'@
$ws.Range("D66").Value = $v
$ws.Rows.Item(66).AutoFit() | Out-Null

# Row 67
$ws.Range("A67").Value = 1339510654663835904
$v = @'
⏩PAKET DESIGN ADOBE ILLUSTRATOR🤩 - Pembuatan Design Menggunakan Adobe Illustrator 18 Module 
https://t.co/V7y8EtJmDu
#PaketDesign #adobeillustrator #pakaimasker #jagajarak #mencucitangan #dirumahaja #kerjadarirumah #belajardirumah #KuliahOnline #IndonesiaMaju #ai #jagakesehatan
'@
$ws.Range("B67").Value = $v
$v = @'
PotekantropusX
'@
$ws.Range("C67").Value = $v
$v = @'
Thu Dec 17 10:00:01 +0000 2020
'@
$ws.Range("D67").Value = $v
$ws.Rows.Item(67).AutoFit() | Out-Null

# Row 68
$ws.Range("A68").Value = 1339224792780534016
$v = @'
Pen belajar, tapi kok raga menolak #kuliahonline
#tugasonline
#rebahan
#uts
'@
$ws.Range("B68").Value = $v
$v = @'
incluede
'@
$ws.Range("C68").Value = $v
$v = @'
Wed Dec 16 15:04:06 +0000 2020
'@
$ws.Range("D68").Value = $v
$ws.Rows.Item(68).AutoFit() | Out-Null

# Row 69
$ws.Range("A69").Value = 1339133165034586112
$v = @'
⏩PAKET DESIGN ADOBE ILLUSTRATOR🤩 - Pembuatan Design Menggunakan Adobe Illustrator 18 Module 
https://t.co/V7y8EtrLeU
#PaketDesign #adobeillustrator #pakaimasker #jagajarak #mencucitangan #dirumahaja #kerjadarirumah #belajardirumah #KuliahOnline #IndonesiaMaju #ai #jagakesehatan
'@
$ws.Range("B69").Value = $v
$v = @'
PotekantropusX
'@
$ws.Range("C69").Value = $v
$v = @'
Wed Dec 16 09:00:01 +0000 2020
'@
$ws.Range("D69").Value = $v
$ws.Rows.Item(69).AutoFit() | Out-Null

# Row 70
$ws.Range("A70").Value = 1339071189486998016
$v = @'
Dear akhir tahun
Tolong jauhkan saya dari dosen yg ngasih tugas seminggu 2 kali, quiz tiap minggu, nanya tiap menjelaskan materi dan presentasi yang kalau ngk jawab di minus, soal ujian tdk sesuai kisi-kisi, dan selalu open cam 😊
#kuliahonline
'@
$ws.Range("B70").Value = $v
$v = @'
ffriskaamalia
'@
$ws.Range("C70").Value = $v
$v = @'
Wed Dec 16 04:53:44 +0000 2020
'@
$ws.Range("D70").Value = $v
$ws.Rows.Item(70).AutoFit() | Out-Null

# Row 71
$ws.Range("A71").Value = 1339052415199244032
$v = @'
Teteplah santuy
Walaupun tugas seperti pasir di pantai...😎
#tugas #kuliahonline
'@
$ws.Range("B71").Value = $v
$v = @'
Ndeanindah
'@
$ws.Range("C71").Value = $v
$v = @'
Wed Dec 16 03:39:08 +0000 2020
'@
$ws.Range("D71").Value = $v
$ws.Rows.Item(71).AutoFit() | Out-Null

# Row 72
$ws.Range("A72").Value = 1338845447058682112
$v = @'
https://t.co/E6oEjMlfxT
Bantuan pelajar dan Mahasiswa
#KuliahOnline 
#indonesiaprokerja
'@
$ws.Range("B72").Value = $v
$v = @'
587bb6eb4c48434
'@
$ws.Range("C72").Value = $v
$v = @'
Tue Dec 15 13:56:43 +0000 2020
'@
$ws.Range("D72").Value = $v
$ws.Rows.Item(72).AutoFit() | Out-Null

# Row 73
$ws.Range("A73").Value = 1338659543665046016
$v = @'
Bacaan dari Blog: Bagaimana ilmu komunikasi lingkungan berkembang dan kenapa ilmu ini penting? https://t.co/b4cinB59AH #blogging #blogger #Bloggers #kuliah #KuliahOnline #komunikasi
'@
$ws.Range("B73").Value = $v
$v = @'
sdpinuji
'@
$ws.Range("C73").Value = $v
$v = @'
Tue Dec 15 01:38:00 +0000 2020
'@
$ws.Range("D73").Value = $v
$ws.Rows.Item(73).AutoFit() | Out-Null

# Row 74
$ws.Range("A74").Value = 1338424929503464960
$v = @'
Ga kerasa kayanya baru kemaren aja ospek tau2 udah UAS aja, semakin kesini semakin nyata kalo hidup itu cuman numpang minum🙂
#UASonline 
#kuliahonline 
#maba
'@
$ws.Range("B74").Value = $v
$v = @'
faiqotuzzahro
'@
$ws.Range("C74").Value = $v
$v = @'
Mon Dec 14 10:05:44 +0000 2020
'@
$ws.Range("D74").Value = $v
$ws.Rows.Item(74).AutoFit() | Out-Null

$ws.Range("J71").Select() | Out-Null

